$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels (row 18) - bold, matching existing label style (e.g. D15/E15)
$ws.Range("D18").Value = "Mean increase"
$ws.Range("D18").Font.Bold = $true
$ws.Range("F18").Value = "Median increase"
$ws.Range("F18").Font.Bold = $true

# New computed values (row 19)
$ws.Range("D19").Formula = "= ((95.488866 / 95.321842) * 100) - 100"
$ws.Range("F19").Formula = "=((95.20038 / 95.22216) * 100) - 100"

# Update selection to match target state
$ws.Range("D20").Select()
